$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above current row 5, shifting rows 5-7 down to 6-8
$ws.Rows("5:5").Insert()

# Copy the date-cell style (stored in column D of the row above) to the new row's D cell
$ws.Range("D4").Copy()
$ws.Range("D5").PasteSpecial(-4122)  # xlPasteFormats

# Fill in the new row's values
$ws.Cells.Item(5, 1).Value = 11
$ws.Cells.Item(5, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(5, 3).Value = "Bíobío"
$ws.Cells.Item(5, 4).Value = 44482
$ws.Cells.Item(5, 5).Value = 8
$ws.Cells.Item(5, 6).Value = 100112022
$ws.Cells.Item(5, 7).Value = "Arveja Verde"
$ws.Cells.Item(5, 8).Value = "Perfection"
$ws.Cells.Item(5, 9).Value = "Primera"
$ws.Cells.Item(5, 10).Value = 130
$ws.Cells.Item(5, 11).Value = 24000
$ws.Cells.Item(5, 12).Value = 25000
$ws.Cells.Item(5, 13).Value = 24385
$ws.Cells.Item(5, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(5, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(5, 16).Value = 975
$ws.Cells.Item(5, 17).Value = 25
$ws.Cells.Item(5, 18).Value = "Hortaliza"
